$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

$ws.Range("B12").Value = 60
$ws.Range("C12").Value = -14
$ws.Range("E12").Value = "46 / 112"
